$wb = $excel.ActiveWorkbook

# --- Update mislabeled text: "Schulden" -> "Rechnungen" aus anderen Abrechnungsperioden ---
$wsSonstige = $wb.Worksheets.Item("Sonstige Ausgaben")
$wsSonstige.Range("B9").Value  = "Rechnungen aus anderen Abrechnungsperioden "
$wsSonstige.Range("B10").Value = "Rechnungen aus anderen Abrechnungsperioden "

# --- Widen column B on "Sonstige Ausgaben" so the longer label fits ---
$wsSonstige.Columns.Item(2).ColumnWidth = 42.17

# --- Move the active selection / active sheet to "Sonstige Ausgaben" ---
$wsSonstige.Activate()
$wsSonstige.Range("B14").Select()
